$d = $word.ActiveDocument

# --- Edit 1: mark "utilizo" as a spell-check flagged word inside the
# "No se utilizo sass y css..." paragraph (paragraph gains proofErr markers
# around "utilizo", matching the existing proofErr pattern already used for
# "sass"/"css"/"tenia"/"comenzo" in the same paragraph). ---
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("No se utilizo ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $find.Parent
    $start = $target.Start
    $end = $target.End
    $oldRange = $d.Range($start, $end)
    $oldRange.Delete()
    $insertionPoint = $d.Range($start, $start)
    $xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">No se </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>utilizo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertionPoint.InsertXML($xml1)
}

# --- Edit 2: append the reviewer's notes as new paragraphs at the very end
# of the document body (before sectPr), following the doc's existing pattern
# of separating paragraphs with a blank paragraph. ---
$endRange = $d.Content
$endRange.Collapse(0)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>No se logra complementar diseño por falta de tiempo.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Se nota que tiene diferencias las columnas y filas en tamaño responsive</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Se deja </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>maquetacion</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Responsive para el final del trabajo.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endRange.InsertXML($xml2)
